$wb = $excel.ActiveWorkbook

# Trade #20 closed at 2026-02-16 22:59:22 - base_strategy UP +0.000%
# The new trade row is appended (as row 21) to both the "All Trades"
# log and the per-strategy "base_strategy" log with identical data.
$sheetNames = @("All Trades", "base_strategy")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Seed row 21 from the last existing row (20) so formatting/cell
    # "shape" (e.g. the blank Exit Price / Exit Reason cells) matches
    # the rest of the trade log, then overwrite the fields that differ
    # for this new trade.
    $ws.Range("A20:Q20").Copy() | Out-Null
    $ws.Range("A21:Q21").PasteSpecial() | Out-Null

    $ws.Cells.Item(21, 1).Value = 20             # Trade #
    # Date (B21) stays "2026-02-16", same as the seeded row - no change needed.
    $ws.Cells.Item(21, 3).Value = "22:59:22"     # Time
    # D: Strategy = base_strategy (unchanged from seed)
    # E: Side = UP (unchanged from seed)
    # F: Entry Price = 0.5 (unchanged from seed)
    # G: Exit Price = blank (unchanged from seed)
    # H: Status = OPEN (unchanged from seed)
    # I: P&L % = 0 (unchanged from seed)
    # J: P&L $ = 0 (unchanged from seed)
    # K: Capital After = 100 (unchanged from seed)
    # L: Entry Slippage (bps) = 0 (unchanged from seed)
    # M: Exit Slippage (bps) = 0 (unchanged from seed)
    # N: Confidence = 0.6 (unchanged from seed)
    # O: Entry Reason = "Normal spread capture: 19600 bps" (unchanged from seed)
    # P: Exit Reason = blank (unchanged from seed)
    # Q: Duration (min) = 0 (unchanged from seed)
}
